$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells before assigning, to preserve exact
# textual representation (e.g. trailing zeros, multi-dot strings) and
# prevent Excel auto-converting these numeric-looking strings to numbers.
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $priceCells) {
  $ws.Range($addr).NumberFormat = "@"
}

# Update Price (column D) values
$ws.Range("D2").Value = "22.464.20"
$ws.Range("D3").Value = "1.566.57"
$ws.Range("D5").Value = "1.001"
$ws.Range("D6").Value = "288.93"
$ws.Range("D7").Value = "0.3674"
$ws.Range("D8").Value = "49.94"
$ws.Range("D9").Value = "0.3374"
$ws.Range("D10").Value = "1.139"
$ws.Range("D11").Value = "0.07490"
$ws.Range("D12").Value = "1.000"
$ws.Range("D13").Value = "21.03"
$ws.Range("D14").Value = "5.988"
$ws.Range("D15").Value = "6.952"
$ws.Range("D16").Value = "1.563.37"
$ws.Range("D17").Value = "0.00001111"
$ws.Range("D18").Value = "90.07"
$ws.Range("D19").Value = "0.06746"
$ws.Range("D21").Value = "6.362"
$ws.Range("D22").Value = "16.20"
$ws.Range("D23").Value = "12.04"
$ws.Range("D24").Value = "22.451.22"
$ws.Range("D25").Value = "2.391"
$ws.Range("D26").Value = "2.628"
$ws.Range("D27").Value = "19.79"
$ws.Range("D28").Value = "149.73"
$ws.Range("D29").Value = "5.021"
$ws.Range("D30").Value = "124.17"
$ws.Range("D31").Value = "1.736.49"
$ws.Range("D32").Value = "1.052"
$ws.Range("D33").Value = "6.153"
$ws.Range("D35").Value = "9.643"
$ws.Range("D36").Value = "0.08291"
$ws.Range("D37").Value = "0.02439"
$ws.Range("D38").Value = "1.330"
$ws.Range("D39").Value = "0.2262"
$ws.Range("D40").Value = "0.06409"
$ws.Range("D41").Value = "5.358"
$ws.Range("D42").Value = "11.20"
$ws.Range("D43").Value = "0.6165"
$ws.Range("D45").Value = "13.77"
$ws.Range("D46").Value = "3.769"
$ws.Range("D47").Value = "0.5782"
$ws.Range("D48").Value = "2.036"
$ws.Range("D49").Value = "125.55"
$ws.Range("D50").Value = "1.225"
$ws.Range("D51").Value = "0.07344"

# Restore default (unstyled) cell style now that the text value is set,
# so the cells keep matching their original (unstyled) appearance.
foreach ($addr in $priceCells) {
  $ws.Range($addr).Style = "Normal"
}

# Update Volume(1h) (column E) values
$ws.Range("E3").Value = "  -0.28%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("E8").Value = "  +0.67%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("E10").Value = "  +0.89%  "
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("E15").Value = "  +1.13%  "
$ws.Range("E17").Value = "  -0.43%  "
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("E20").Value = "  -0.20%  "
$ws.Range("E21").Value = "  +3.07%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("E23").Value = "  +1.49%  "
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("E25").Value = "  +1.01%  "
$ws.Range("E26").Value = "  +2.96%  "
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("E28").Value = "  +1.82%  "
$ws.Range("E29").Value = "  +0.42%  "
$ws.Range("E30").Value = "  -0.43%  "
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("E32").Value = "  +5.99%  "
$ws.Range("E34").Value = "  +2.58%  "
$ws.Range("E35").Value = "  -1.07%  "
$ws.Range("E36").Value = "  -1.63%  "
$ws.Range("E37").Value = "  -0.46%  "
$ws.Range("E38").Value = "  -3.73%  "
$ws.Range("E39").Value = "  +0.40%  "
$ws.Range("E40").Value = "  -1.10%  "
$ws.Range("E41").Value = "  -0.69%  "
$ws.Range("E42").Value = "  -1.22%  "
$ws.Range("E43").Value = "  -1.02%  "
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("E45").Value = "  -1.43%  "
$ws.Range("E46").Value = "  -1.14%  "
$ws.Range("E47").Value = "  -0.18%  "
$ws.Range("E48").Value = "  -0.95%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("E50").Value = "  -0.20%  "
$ws.Range("E51").Value = "  +0.55%  "
